{"js": "// Update the date line and the 25 division-problem cells to the new\n// worksheet values (see commit: \"Update master to output generated at\n// 4250d90\"). Every old value in this document is unique, so a\n// search-and-replace keyed on the exact old text is unambiguous.\nconst replacements = [\n  [\"2024-08-12 Monday\", \"2024-08-13 Tuesday\"],\n  [\"849\u00f78=106, 1\", \"831\u00f75=166, 1\"],\n  [\"144\u00f76=24, 0\", \"555\u00f74=138, 3\"],\n  [\"736\u00f74=184, 0\", \"230\u00f79=25, 5\"],\n  [\"948\u00f72=474, 0\", \"461\u00f73=153, 2\"],\n  [\"963\u00f73=321, 0\", \"551\u00f74=137, 3\"],\n  [\"926\u00f76=154, 2\", \"547\u00f76=91, 1\"],\n  [\"357\u00f77=51, 0\", \"942\u00f77=134, 4\"],\n  [\"536\u00f72=268, 0\", \"816\u00f76=136, 0\"],\n  [\"153\u00f74=38, 1\", \"576\u00f78=72, 0\"],\n  [\"280\u00f77=40, 0\", \"181\u00f73=60, 1\"],\n  [\"287\u00f79=31, 8\", \"788\u00f79=87, 5\"],\n  [\"639\u00f73=213, 0\", \"381\u00f72=190, 1\"],\n  [\"426\u00f76=71, 0\", \"719\u00f72=359, 1\"],\n  [\"142\u00f78=17, 6\", \"416\u00f74=104, 0\"],\n  [\"560\u00f73=186, 2\", \"708\u00f79=78, 6\"],\n  [\"929\u00f77=132, 5\", \"161\u00f77=23, 0\"],\n  [\"283\u00f74=70, 3\", \"535\u00f73=178, 1\"],\n  [\"908\u00f72=454, 0\", \"943\u00f79=104, 7\"],\n  [\"689\u00f79=76, 5\", \"159\u00f77=22, 5\"],\n  [\"871\u00f78=108, 7\", \"443\u00f78=55, 3\"],\n  [\"367\u00f78=45, 7\", \"838\u00f78=104, 6\"],\n  [\"384\u00f76=64, 0\", \"794\u00f79=88, 2\"],\n  [\"757\u00f74=189, 1\", \"698\u00f72=349, 0\"],\n  [\"318\u00f75=63, 3\", \"725\u00f73=241, 2\"],\n  [\"500\u00f75=100, 0\", \"153\u00f77=21, 6\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and the 25 division-problem cells to the new\n# worksheet values (see commit: \"Update master to output generated at\n# 4250d90\"). Every old value in this document is unique, so a\n# find/replace keyed on the exact old text is unambiguous.\n\n$replacements = @(\n    @(\"2024-08-12 Monday\", \"2024-08-13 Tuesday\"),\n    @(\"849\u00f78=106, 1\", \"831\u00f75=166, 1\"),\n    @(\"144\u00f76=24, 0\", \"555\u00f74=138, 3\"),\n    @(\"736\u00f74=184, 0\", \"230\u00f79=25, 5\"),\n    @(\"948\u00f72=474, 0\", \"461\u00f73=153, 2\"),\n    @(\"963\u00f73=321, 0\", \"551\u00f74=137, 3\"),\n    @(\"926\u00f76=154, 2\", \"547\u00f76=91, 1\"),\n    @(\"357\u00f77=51, 0\", \"942\u00f77=134, 4\"),\n    @(\"536\u00f72=268, 0\", \"816\u00f76=136, 0\"),\n    @(\"153\u00f74=38, 1\", \"576\u00f78=72, 0\"),\n    @(\"280\u00f77=40, 0\", \"181\u00f73=60, 1\"),\n    @(\"287\u00f79=31, 8\", \"788\u00f79=87, 5\"),\n    @(\"639\u00f73=213, 0\", \"381\u00f72=190, 1\"),\n    @(\"426\u00f76=71, 0\", \"719\u00f72=359, 1\"),\n    @(\"142\u00f78=17, 6\", \"416\u00f74=104, 0\"),\n    @(\"560\u00f73=186, 2\", \"708\u00f79=78, 6\"),\n    @(\"929\u00f77=132, 5\", \"161\u00f77=23, 0\"),\n    @(\"283\u00f74=70, 3\", \"535\u00f73=178, 1\"),\n    @(\"908\u00f72=454, 0\", \"943\u00f79=104, 7\"),\n    @(\"689\u00f79=76, 5\", \"159\u00f77=22, 5\"),\n    @(\"871\u00f78=108, 7\", \"443\u00f78=55, 3\"),\n    @(\"367\u00f78=45, 7\", \"838\u00f78=104, 6\"),\n    @(\"384\u00f76=64, 0\", \"794\u00f79=88, 2\"),\n    @(\"757\u00f74=189, 1\", \"698\u00f72=349, 0\"),\n    @(\"318\u00f75=63, 3\", \"725\u00f73=241, 2\"),\n    @(\"500\u00f75=100, 0\", \"153\u00f77=21, 6\")\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $r = $d.Content\n    $find = $r.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1          # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        Write-Host \"NOT FOUND: $oldText\"\n    }\n}\n"}
